$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds for existing rows 4 and 5
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6

$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63

# Delete row 6 (EGYPT - PREMIER LEAGUE / ZED vs Al Ahly) entirely
$ws.Rows(6).Delete()

# After deletion, old row 7 (West Brom vs Burnley) is now row 6.
# Update its odds values that changed
$ws.Range("G6").Value = 2.45
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 3
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.62
$ws.Range("W6").Value = 7
$ws.Range("AI6").Value = 13
$ws.Range("AS6").Value = 251
$ws.Range("AW6").Value = 4.75
$ws.Range("AZ6").Value = 51
$ws.Range("BA6").Value = 81

# Delete the old row 8 (Poland - Division 1), which after the first deletion is now row 7
$ws.Rows(7).Delete()
